$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "243.08"
Set-TextValue "D3" "22.98"
Set-TextValue "D4" "5.393"
Set-TextValue "D5" "0.05909"
Set-TextValue "D6" "3.455"
Set-TextValue "D7" "6.557"
Set-TextValue "D8" "0.8144"
Set-TextValue "D9" "0.9153"
Set-TextValue "D10" "0.1418"
Set-TextValue "D11" "0.07421"
Set-TextValue "D12" "0.03273"
Set-TextValue "D13" "0.03062"
Set-TextValue "D14" "0.09342"
Set-TextValue "D15" "3.846"
Set-TextValue "D16" "0.001563"
Set-TextValue "D17" "0.04674"
Set-TextValue "D18" "0.0005911"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue "D19" "0.005897"
Set-TextValue "D20" "0.001294"
Set-TextValue "D21" "0.004909"
Set-TextValue "D22" "0.00009501"
Set-TextValue "D24" "2.151"
Set-TextValue "D40" "0.03954"
Set-TextValue "D41" "0.006177"
Set-TextValue "D42" "0.1071"
Set-TextValue "D43" "0.002700"
Set-TextValue "D44" "0.008118"
Set-TextValue "D45" "0.00005197"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
